# "bank data is stopped and calendar is started"
#
# Finishes row 6 (2nd sub-step of task 2), then appends three new tasks:
#   3) رسیدن به مرحله استخراج اکسل هر کدوم از بانک ها        - row 7
#   4) ساخت اکسل کامل                                         - rows 8-10 (merged ردیف/شرح وظایف)
#   5) ساخت تقویم                                              - row 11
# and turns the plain "=F-E" formulas into one shared formula over G2:G11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Copy-CellFormat([string]$fromAddr, [string]$toAddr) {
    $ws.Range($fromAddr).Copy()
    $ws.Range($toAddr).PasteSpecial($xlPasteFormats)
}

function Set-TextValue([string]$addr, [string]$text, [string]$likeAddr) {
    # Force literal text (column D holds Jalali-calendar date *strings*, not
    # real Excel date serials) and then restore the plain/general formatting
    # that the rest of the column already uses.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    Copy-CellFormat $likeAddr $addr
}

# ---- finish row 6 (second sub-step of task 2) ----
$ws.Range("F6").Value = 2
$ws.Range("G6").Formula = "=F6-E6"

# ---- task 3: single row ----
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "رسیدن به مرحله استخراج اکسل هر کدوم از بانک ها"
$ws.Range("C7").Value = "متاسفانه تا این مرحله یادم رفت ساعت کاری رو درج کنم و این یک زمان حدودی انجام کار است که در این تاریخ زده می شود"
Set-TextValue "D7" "1403/09/11" "D2"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 15
$ws.Range("G7").Formula = "=F7-E7"

Copy-CellFormat "A5" "A7"
Copy-CellFormat "C5" "B7"
Copy-CellFormat "C5" "C7"

# ---- task 4: spans rows 8-10 ----
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "ساخت اکسل کامل "
$ws.Range("C8").Value = "در این مرحله یک اکسل بزرگ از اطلاعات همه ی بانک ها ساخته می شود"
Set-TextValue "D8" "1403/09/11" "D2"
$ws.Range("E8").Value = 19
$ws.Range("F8").Value = 21.5
$ws.Range("G8").Formula = "=F8-E8"

$ws.Range("C9").Value = "در انجام این مرحله یک سری اکسل امکان خوانده شدن نداشتن که این مدت به رفع این مشکل گذشت"
Set-TextValue "D9" "1403/09/13" "D2"
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 2
$ws.Range("G9").Formula = "=F9-E9"

Set-TextValue "D10" "1403/09/13" "D2"
$ws.Range("E10").Value = 13.5
$ws.Range("F10").Value = 14.5
$ws.Range("G10").Formula = "=F10-E10"

Copy-CellFormat "A2" "A8"
Copy-CellFormat "A2" "A9"
Copy-CellFormat "A2" "A10"
Copy-CellFormat "B2" "B8"
Copy-CellFormat "B2" "B9"
Copy-CellFormat "B2" "B10"
Copy-CellFormat "D2" "C8"
Copy-CellFormat "C2" "C9"
Copy-CellFormat "C2" "C10"

$ws.Range("A8:A10").Merge()
$ws.Range("B8:B10").Merge()
$ws.Range("C9:C10").Merge()

# ---- task 5: single row ----
$ws.Range("A11").Value = 5
$ws.Range("B11").Value = "ساخت تقویم "
$ws.Range("C11").Value = "در این مرحله دیتای سال ۱۴۰۳ تا حد شکل گیری تقویم جمع َآوری شد و تقویم ساخته شد ولی اطلاعات ماه قمری جهت تعیطلات اشتباه است"
Set-TextValue "D11" "1403/09/14" "D2"
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 3
$ws.Range("G11").Formula = "=F11-E11"

Copy-CellFormat "A5" "A11"
Copy-CellFormat "D2" "B11"
Copy-CellFormat "C5" "C11"

# ---- row heights to match the two-line wrapped descriptions ----
$ws.Rows.Item(7).RowHeight = 28.8
$ws.Rows.Item(11).RowHeight = 28.8

$ws.Range("C11").Select()
